$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new rows, copying formatting from the existing data row (row 16) ---
# Row 16 currently holds the single "1901" period line; Excel duplicates it
# (copy + insert) to grow the table to 5 period rows, matching how a user would
# add new Estado de Cuenta period lines above the closing (last) row of the block.
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows.Item(16).Copy() | Out-Null
    $ws.Rows.Item(17).Insert() | Out-Null
}

# --- Fill in the period / value data for the five rows ---
# New periods (most recent first): 2006, 2005, 2004, 2003, 1901
$ws.Range("E16").Value = "2006"
$ws.Range("F16").Value = 35112
$ws.Range("G16").Value = 877803

$ws.Range("E17").Value = "2005"
$ws.Range("F17").Value = 35112
$ws.Range("G17").Value = 877803

$ws.Range("E18").Value = "2004"
$ws.Range("F18").Value = 35112
$ws.Range("G18").Value = 877803

$ws.Range("E19").Value = "2003"
$ws.Range("F19").Value = 35112
$ws.Range("G19").Value = 877803

$ws.Range("E20").Value = "1901"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 877803

# --- Update the summary fields above the table ---
# Valor Mora total = sum of the five period values
$ws.Range("E11").Value = 171697
# Cant. Periodos = number of period rows now in the table
$ws.Range("F13").Value = 5
